# Auto-generated Excel COM-interop script
# Applies the cryptos.xlsx data refresh described in the commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.605.10"
$ws.Range("E2").Value = "  -1.17%  "
$ws.Range("D3").Value = "2.534.27"
$ws.Range("E3").Value = "  -1.43%  "
$__style = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "308.35"
$ws.Range("D5").Style = $__style
$ws.Range("E5").Value = "  -2.10%  "
$__style = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "100.75"
$ws.Range("D6").Style = $__style
$ws.Range("E6").Value = "  +4.23%  "
$ws.Range("E7").Value = "  -1.18%  "
$ws.Range("E8").Value = "  +0.09%  "
$__style = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.528"
$ws.Range("D9").Style = $__style
$ws.Range("E9").Value = "  -2.04%  "
$__style = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.88"
$ws.Range("D10").Style = $__style
$ws.Range("E10").Value = "  +1.19%  "
$ws.Range("E11").Value = "  -1.14%  "
$__style = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.35"
$ws.Range("D12").Style = $__style
$ws.Range("E12").Value = "  -1.24%  "
$ws.Range("E13").Value = "  +0.04%  "
$ws.Range("D14").Value = "2.934.51"
$ws.Range("E14").Value = "  -1.00%  "
$__style = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.90"
$ws.Range("D15").Style = $__style
$ws.Range("E15").Value = "  +5.65%  "
$ws.Range("D16").Value = "2.492.26"
$ws.Range("E16").Value = "  -4.02%  "
$__style = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.813"
$ws.Range("D17").Style = $__style
$ws.Range("E17").Value = "  -3.52%  "
$ws.Range("D18").Value = "42.611.12"
$ws.Range("E18").Value = "  -1.20%  "
$ws.Range("E19").Value = "  -0.97%  "
$ws.Range("D20").Value = "0.0₃0953"
$ws.Range("E20").Value = "  -0.81%  "
$__style = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.22"
$ws.Range("D21").Style = $__style
$ws.Range("E21").Value = "  -2.98%  "
$__style = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "69.45"
$ws.Range("D22").Style = $__style
$ws.Range("E22").Value = "  +0.15%  "
$__style = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "243.54"
$ws.Range("D23").Style = $__style
$ws.Range("E23").Value = "  -3.79%  "
$ws.Range("E24").Value = "  -2.01%  "
$ws.Range("E25").Value = "  -1.23%  "
$ws.Range("E26").Value = "  +0.07%  "
$__style = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "26.00"
$ws.Range("D27").Style = $__style
$ws.Range("E27").Value = "  -3.33%  "
$ws.Range("E28").Value = "  -3.85%  "
$__style = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "39.27"
$ws.Range("D29").Style = $__style
$ws.Range("E29").Value = "  -2.38%  "
$__style = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "10.15"
$ws.Range("D30").Style = $__style
$ws.Range("E30").Value = "  -1.20%  "
$__style = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.79"
$ws.Range("D31").Style = $__style
$ws.Range("E31").Value = "  -0.72%  "
$__style = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "155.57"
$ws.Range("D32").Style = $__style
$ws.Range("E32").Value = "  +0.69%  "
$ws.Range("E33").Value = "  +11.74%  "
$ws.Range("E34").Value = "  -1.54%  "
$ws.Range("E35").Value = "  -2.83%  "
$__style = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.35"
$ws.Range("D36").Style = $__style
$ws.Range("E36").Value = "  -3.16%  "
$ws.Range("E37").Value = "  -4.76%  "
$__style = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.16"
$ws.Range("D38").Style = $__style
$ws.Range("E38").Value = "  -6.43%  "
$__style = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.113"
$ws.Range("D39").Style = $__style
$ws.Range("E39").Value = "  +0.98%  "
$ws.Range("E40").Value = "  +0.48%  "
$__style = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.35"
$ws.Range("D41").Style = $__style
$ws.Range("E41").Value = "  +10.05%  "
$__style = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "22.13"
$ws.Range("D42").Style = $__style
$ws.Range("E42").Value = "  -1.16%  "
$ws.Range("E43").Value = "  +0.08%  "
$ws.Range("E44").Value = "  +1.72%  "
$ws.Range("E45").Value = "  -1.87%  "
$ws.Range("D46").Value = "1.973.57"
$ws.Range("E46").Value = "  -1.18%  "
$ws.Range("E47").Value = "  -0.15%  "
$__style = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "81.05"
$ws.Range("D48").Style = $__style
$ws.Range("E48").Value = "  -2.26%  "
$ws.Range("B49").Value = "SEI"
$ws.Range("C49").Value = "https://coinranking.com/coin/8nxCqs-uj+sei-sei"
$__style = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.857"
$ws.Range("D49").Style = $__style
$ws.Range("E49").Value = "  +10.60%  "
$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").Value = "2.731.00"
$ws.Range("E50").Value = "  -2.99%  "
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$__style = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.192"
$ws.Range("D51").Style = $__style
$ws.Range("E51").Value = "  -0.89%  "
